$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Dario Granizo",    "1719534925", "GOL201900175", "javierdar5@yopmail.com"),
    @("Alvaro Garcia",    "0951672849", "201900015",    "alvarogm2009@yopmail.com"),
    @("Joan Achi",        "0924131766", "",             "joandaniellaachimantilla@yopmail.com"),
    @("Jordan Acaro",     "0941571069", "201800175",    "jordanacaro64@yopmail.com"),
    @("Gavriela Aguilar", "0925624934", "",             "gabrielaaguilar97@yopmail.com"),
    @("Diana Calderon",   "0104537980", "201800157",    "gaby.14x@yopmail.com")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    if ($rec[2] -ne "") {
        $ws.Cells.Item($row, 3).Value = $rec[2]
    }
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row = $row + 1
}

[void]$ws.Range("B14").Select()
